$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap "backup@backdoor.com, System" -> "System, backup@backdoor.com" ---
# (keeps any trailing ", system" suffix intact) across every G-column cell
# that currently holds that text.
$gCells = @("G2","G4","G5","G8","G29","G31","G32","G35","G56","G58","G59","G62", `
            "G83","G84","G85","G109","G110","G111","G135","G136","G137")

foreach ($ref in $gCells) {
    $cell = $ws.Range($ref)
    $old = $cell.Value2
    if ($old -ne $null) {
        $new = $old -replace "backup@backdoor\.com, System", "System, backup@backdoor.com"
        $cell.Value = $new
    }
}

# --- 2. Numeric tweaks ---
$ws.Range("L7").Value = 3
$ws.Range("L8").Value = 21

$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 4
$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 4
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 4

# --- 3. Column I (9) a little wider (stored width 10 -> 14) ---
$ws.Columns.Item(9).ColumnWidth = 13.17

# --- 4. Rows 104 / 130 / 156 move from "Pending" (yellow) to "Not Recorded" (pink) ---
$headerRows = @(104,130,156)
foreach ($r in $headerRows) {
    $rng = $ws.Range("A" + $r + ":I" + $r)
    $rng.Interior.Color = 12695295   # RGB(255,182,193) "light pink" fill
    $rng.Font.Color = 0              # black text
    $ws.Range("I" + $r).Value = "Not Recorded"
}
